$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Consolidate the title text run(s) on the slide itself into a single run.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Delete()
$title.Text = "Here is a single header"

# Consolidate the speaker-notes text run(s) into a single run.
$notes = $s.NotesPage
$notesText = $notes.Shapes.Item(2).TextFrame.TextRange
$notesText.Text = "and here are some notes"
